$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet gains a new "数据结构" (data-structure) table at the very top,
# pushing all of the existing content down by 39 rows (with a further gap
# before the new final summary row). Replicate that by inserting 39 blank
# rows at the top, then populating the newly created rows.
$ws.Rows("1:39").Insert()

# --- new header block -------------------------------------------------
$ws.Range("A1").Value = "数据结构"

$ws.Range("A2").Value = "Header"
$ws.Range("B2").Value = "包类型-SID"
$ws.Range("C2").Value = "序列号"
$ws.Range("D2").Value = "DataLength"
$ws.Range("E2").Value = "Data"
$ws.Range("F2").Value = "CRC16"

$ws.Range("A3").Value = "1Byte"
$ws.Range("B3").Value = "4bit-4bit"
$ws.Range("C3").Value = "1Byte"
$ws.Range("D3").Value = "1Byte"
$ws.Range("E3").Value = "nByte"
$ws.Range("F3").Value = "2Byte"

# Rows 4/5 only carry formatting (vertically centred, no wrap) on column B,
# no text value.
$ws.Range("B4").VerticalAlignment = -4108
$ws.Range("B4").WrapText = $false
$ws.Range("B5").VerticalAlignment = -4108
$ws.Range("B5").WrapText = $false

$ws.Range("A6").Value = "包类型"

$ws.Range("A7").Value = "SYN"
$ws.Range("B7").Value = "bit7"

$ws.Range("A8").Value = "ACK"
$ws.Range("B8").Value = "bit6"

$ws.Range("A9").Value = "SID"
$ws.Range("B9").Value = "bit0~bit3"

# --- new trailing note under the existing content ----------------------
$ws.Range("A78").Value = "通讯为异步"

# --- column / view tweaks ----------------------------------------------
$ws.Columns("A").ColumnWidth = 11.08984375

# Scroll back to the top and leave the selection on B9, matching the saved
# view state of the edited workbook.
$ws.Range("B9").Select()

Write-Host "applied protocol header table + async note"
